$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 343,
# shifting every subsequent record (old rows 343-418) down by one row
# (to 344-419). Insert a blank row at 343 first -- Excel's native Insert
# shifts the data down and carries the row-above's formatting (including
# the date number format on column D) onto the new row.
$ws.Rows("343:343").Insert()

# Populate the newly inserted row 343 with the new record's data.
$ws.Range("A343").Value = 8
$ws.Range("B343").Value = "Terminal La Palmera de La Serena"
$ws.Range("C343").Value = "Coquimbo"
$ws.Range("D343").Value = (Get-Date -Year 2022 -Month 5 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E343").Value = 4
$ws.Range("F343").Value = 100114001
$ws.Range("G343").Value = "Papa"
$ws.Range("H343").Value = "Asterix"
$ws.Range("I343").Value = "1a (cosecha)"
$ws.Range("J343").Value = 2520
$ws.Range("K343").Value = 8000
$ws.Range("L343").Value = 9000
$ws.Range("M343").Value = 8500
$ws.Range("N343").Value = "`$/saco 25 kilos"
$ws.Range("O343").Value = "Región de Los Lagos"
$ws.Range("P343").Value = 340
$ws.Range("Q343").Value = 25
$ws.Range("R343").Value = "Hortaliza"
